# Updated cryptos list on Fri Apr 14 02:52:12 UTC 2023 with GitHub Actions.
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns for every coin row,
# and swaps the BinanceUSD / ShibaInu rows (17 <-> 18) to match the new
# ranking order/data pulled from the coinranking API.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is the new cell text. Numeric-looking price strings (e.g.
# "1.003", "0.5200") are prefixed with a leading apostrophe so Excel keeps
# storing them as text - same as the source data - instead of silently
# coercing them to a Number and dropping significant trailing zeros.
$updates = @(
    @{ Cell = 'D2'; Value = '30.686.47' },
    @{ Cell = 'E2'; Value = '  +1.90%  ' },
    @{ Cell = 'D3'; Value = '2.112.22' },
    @{ Cell = 'E3'; Value = '  +10.81%  ' },
    @{ Cell = 'D4'; Value = '''1.003' },
    @{ Cell = 'E4'; Value = '  +0.28%  ' },
    @{ Cell = 'D5'; Value = '''331.84' },
    @{ Cell = 'E5'; Value = '  +3.84%  ' },
    @{ Cell = 'D6'; Value = '''1.002' },
    @{ Cell = 'E6'; Value = '  +0.28%  ' },
    @{ Cell = 'D7'; Value = '''0.5200' },
    @{ Cell = 'E7'; Value = '  +2.97%  ' },
    @{ Cell = 'D8'; Value = '''0.4364' },
    @{ Cell = 'E8'; Value = '  +7.21%  ' },
    @{ Cell = 'D9'; Value = '''0.08963' },
    @{ Cell = 'E9'; Value = '  +7.51%  ' },
    @{ Cell = 'D10'; Value = '''45.13' },
    @{ Cell = 'E10'; Value = '  +6.78%  ' },
    @{ Cell = 'D11'; Value = '''1.171' },
    @{ Cell = 'E11'; Value = '  +6.18%  ' },
    @{ Cell = 'D12'; Value = '''24.81' },
    @{ Cell = 'E12'; Value = '  +3.83%  ' },
    @{ Cell = 'D13'; Value = '2.126.90' },
    @{ Cell = 'E13'; Value = '  +11.93%  ' },
    @{ Cell = 'D14'; Value = '''6.778' },
    @{ Cell = 'E14'; Value = '  +6.11%  ' },
    @{ Cell = 'D15'; Value = '''7.668' },
    @{ Cell = 'E15'; Value = '  +6.27%  ' },
    @{ Cell = 'D16'; Value = '''97.42' },
    @{ Cell = 'E16'; Value = '  +5.51%  ' },
    @{ Cell = 'B17'; Value = 'BinanceUSD' },
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd' },
    @{ Cell = 'D17'; Value = '''1.002' },
    @{ Cell = 'E17'; Value = '  +0.32%  ' },
    @{ Cell = 'B18'; Value = 'ShibaInu' },
    @{ Cell = 'C18'; Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib' },
    @{ Cell = 'D18'; Value = '''0.00001136' },
    @{ Cell = 'E18'; Value = '  +3.80%  ' },
    @{ Cell = 'D19'; Value = '''0.06611' },
    @{ Cell = 'E19'; Value = '  +1.84%  ' },
    @{ Cell = 'D20'; Value = '''19.15' },
    @{ Cell = 'E20'; Value = '  +4.90%  ' },
    @{ Cell = 'D21'; Value = '''6.418' },
    @{ Cell = 'E21'; Value = '  +8.38%  ' },
    @{ Cell = 'D22'; Value = '''1.000' },
    @{ Cell = 'E22'; Value = '  +0.11%  ' },
    @{ Cell = 'D23'; Value = '30.876.77' },
    @{ Cell = 'E23'; Value = '  +2.51%  ' },
    @{ Cell = 'D24'; Value = '''11.97' },
    @{ Cell = 'E24'; Value = '  +5.66%  ' },
    @{ Cell = 'D25'; Value = '2.365.31' },
    @{ Cell = 'E25'; Value = '  +11.40%  ' },
    @{ Cell = 'D26'; Value = '''2.262' },
    @{ Cell = 'E26'; Value = '  +3.39%  ' },
    @{ Cell = 'D27'; Value = '''22.94' },
    @{ Cell = 'E27'; Value = '  +5.61%  ' },
    @{ Cell = 'D28'; Value = '''2.548' },
    @{ Cell = 'E28'; Value = '  +11.37%  ' },
    @{ Cell = 'D29'; Value = '''163.38' },
    @{ Cell = 'E29'; Value = '  +0.51%  ' },
    @{ Cell = 'D30'; Value = '''133.98' },
    @{ Cell = 'E30'; Value = '  +4.32%  ' },
    @{ Cell = 'D31'; Value = '''1.180' },
    @{ Cell = 'E31'; Value = '  +3.63%  ' },
    @{ Cell = 'D32'; Value = '''0.1069' },
    @{ Cell = 'E32'; Value = '  +2.66%  ' },
    @{ Cell = 'D33'; Value = '''6.197' },
    @{ Cell = 'E33'; Value = '  +4.18%  ' },
    @{ Cell = 'D34'; Value = '''3.915' },
    @{ Cell = 'E34'; Value = '  +3.48%  ' },
    @{ Cell = 'D35'; Value = '''1.519' },
    @{ Cell = 'E35'; Value = '  +27.64%  ' },
    @{ Cell = 'D36'; Value = '''0.02579' },
    @{ Cell = 'E36'; Value = '  +5.10%  ' },
    @{ Cell = 'D37'; Value = '''5.607' },
    @{ Cell = 'E37'; Value = '  +4.68%  ' },
    @{ Cell = 'D38'; Value = '''0.06774' },
    @{ Cell = 'E38'; Value = '  +6.30%  ' },
    @{ Cell = 'D39'; Value = '''9.525' },
    @{ Cell = 'E39'; Value = '  +10.94%  ' },
    @{ Cell = 'D40'; Value = '''12.67' },
    @{ Cell = 'E40'; Value = '  +11.40%  ' },
    @{ Cell = 'D41'; Value = '''0.2242' },
    @{ Cell = 'E41'; Value = '  +4.65%  ' },
    @{ Cell = 'D42'; Value = '''0.6778' },
    @{ Cell = 'E42'; Value = '  +3.88%  ' },
    @{ Cell = 'D43'; Value = '''1.250' },
    @{ Cell = 'E43'; Value = '  +3.02%  ' },
    @{ Cell = 'D44'; Value = '''14.30' },
    @{ Cell = 'E44'; Value = '  +6.81%  ' },
    @{ Cell = 'D45'; Value = '''1.001' },
    @{ Cell = 'E45'; Value = '  +0.19%  ' },
    @{ Cell = 'D46'; Value = '''0.6288' },
    @{ Cell = 'E46'; Value = '  +3.53%  ' },
    @{ Cell = 'D47'; Value = '''2.245' },
    @{ Cell = 'E47'; Value = '  +1.96%  ' },
    @{ Cell = 'D48'; Value = '''3.661' },
    @{ Cell = 'E48'; Value = '  +1.35%  ' },
    @{ Cell = 'D49'; Value = '''1.272' },
    @{ Cell = 'E49'; Value = '  +5.20%  ' },
    @{ Cell = 'D50'; Value = '''126.94' },
    @{ Cell = 'E50'; Value = '  +4.57%  ' },
    @{ Cell = 'D51'; Value = '''83.14' },
    @{ Cell = 'E51'; Value = '  +5.43%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
